$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing column D values (rows 2-67)
$ws.Cells.Item(2, 4).Value = 38.14363599951554
$ws.Cells.Item(3, 4).Value = 15.45791201089473
$ws.Cells.Item(4, 4).Value = 54.6262413454148
$ws.Cells.Item(5, 4).Value = 35.98424184149363
$ws.Cells.Item(6, 4).Value = 32.95278413800352
$ws.Cells.Item(7, 4).Value = 15.1995721316121
$ws.Cells.Item(8, 4).Value = 35.84458650031529
$ws.Cells.Item(9, 4).Value = 14.19217523338042
$ws.Cells.Item(10, 4).Value = 49.32947002795839
$ws.Cells.Item(11, 4).Value = 30.3098107726709
$ws.Cells.Item(12, 4).Value = 32.35899848396569
$ws.Cells.Item(13, 4).Value = 15.62068213575621
$ws.Cells.Item(14, 4).Value = 35.89243503191006
$ws.Cells.Item(15, 4).Value = 12.88367483463914
$ws.Cells.Item(16, 4).Value = 59.01210110543519
$ws.Cells.Item(17, 4).Value = 36.66593134392486
$ws.Cells.Item(18, 4).Value = 33.5295730065754
$ws.Cells.Item(19, 4).Value = 15.99774837246488
$ws.Cells.Item(20, 4).Value = 33.16950160844431
$ws.Cells.Item(21, 4).Value = 10.80300732285277
$ws.Cells.Item(22, 4).Value = 48.58641969174509
$ws.Cells.Item(23, 4).Value = 27.73229493325686
$ws.Cells.Item(24, 4).Value = 30.47536094911921
$ws.Cells.Item(25, 4).Value = 14.07420301758494
$ws.Cells.Item(26, 4).Value = 26.30646878231453
$ws.Cells.Item(27, 4).Value = 8.671687335283092
$ws.Cells.Item(28, 4).Value = 47.80494953793603
$ws.Cells.Item(29, 4).Value = 28.70320527815434
$ws.Cells.Item(30, 4).Value = 28.81466040126037
$ws.Cells.Item(31, 4).Value = 12.90716267834566
$ws.Cells.Item(32, 4).Value = 24.807186663795
$ws.Cells.Item(33, 4).Value = 8.003898898818537
$ws.Cells.Item(34, 4).Value = 47.0814702014259
$ws.Cells.Item(35, 4).Value = 26.13566830483837
$ws.Cells.Item(36, 4).Value = 27.80768597275787
$ws.Cells.Item(37, 4).Value = 13.01520288643518
$ws.Cells.Item(38, 4).Value = 28.60001972228088
$ws.Cells.Item(39, 4).Value = 8.00283514125133
$ws.Cells.Item(40, 4).Value = 35.59701910003867
$ws.Cells.Item(41, 4).Value = 18.79637894593471
$ws.Cells.Item(42, 4).Value = 21.03089170557897
$ws.Cells.Item(43, 4).Value = 8.146251806105521
$ws.Cells.Item(44, 4).Value = 24.39553892066911
$ws.Cells.Item(45, 4).Value = 7.368493146925243
$ws.Cells.Item(46, 4).Value = 38.60029021358826
$ws.Cells.Item(47, 4).Value = 21.46686139770912
$ws.Cells.Item(48, 4).Value = 18.25885883414072
$ws.Cells.Item(49, 4).Value = 6.268733985477996
$ws.Cells.Item(50, 4).Value = 24.22079353443872
$ws.Cells.Item(51, 4).Value = 6.950013272826726
$ws.Cells.Item(52, 4).Value = 44.72056929447555
$ws.Cells.Item(53, 4).Value = 27.39625794597888
$ws.Cells.Item(54, 4).Value = 19.71800804442211
$ws.Cells.Item(55, 4).Value = 7.929880485607629
$ws.Cells.Item(56, 4).Value = 22.30064586507469
$ws.Cells.Item(57, 4).Value = 6.764499708051605
$ws.Cells.Item(58, 4).Value = 44.98516009364589
$ws.Cells.Item(59, 4).Value = 28.27066293837953
$ws.Cells.Item(60, 4).Value = 21.02274696651592
$ws.Cells.Item(61, 4).Value = 8.502292553170685
$ws.Cells.Item(62, 4).Value = 20.29624540007003
$ws.Cells.Item(63, 4).Value = 5.754442569870586
$ws.Cells.Item(64, 4).Value = 42.46372755399076
$ws.Cells.Item(65, 4).Value = 25.78912921082178
$ws.Cells.Item(66, 4).Value = 20.18724028989515
$ws.Cells.Item(67, 4).Value = 8.226707483836581

# Append new rows 68-109
$ws.Cells.Item(68, 1).Value = 2018
$ws.Cells.Item(68, 2).Value = "Costa"
$ws.Cells.Item(68, 3).Value = "Pobreza"
$ws.Cells.Item(68, 4).Value = 23.44433224822351
$ws.Cells.Item(69, 1).Value = 2018
$ws.Cells.Item(69, 2).Value = "Costa"
$ws.Cells.Item(69, 3).Value = "Pobreza extrema"
$ws.Cells.Item(69, 4).Value = 7.106656413560558
$ws.Cells.Item(70, 1).Value = 2018
$ws.Cells.Item(70, 2).Value = "Oriente"
$ws.Cells.Item(70, 3).Value = "Pobreza"
$ws.Cells.Item(70, 4).Value = 45.47921405004109
$ws.Cells.Item(71, 1).Value = 2018
$ws.Cells.Item(71, 2).Value = "Oriente"
$ws.Cells.Item(71, 3).Value = "Pobreza extrema"
$ws.Cells.Item(71, 4).Value = 26.33234339960892
$ws.Cells.Item(72, 1).Value = 2018
$ws.Cells.Item(72, 2).Value = "Sierra"
$ws.Cells.Item(72, 3).Value = "Pobreza"
$ws.Cells.Item(72, 4).Value = 20.19220508237685
$ws.Cells.Item(73, 1).Value = 2018
$ws.Cells.Item(73, 2).Value = "Sierra"
$ws.Cells.Item(73, 3).Value = "Pobreza extrema"
$ws.Cells.Item(73, 4).Value = 7.631384074461965
$ws.Cells.Item(74, 1).Value = 2019
$ws.Cells.Item(74, 2).Value = "Costa"
$ws.Cells.Item(74, 3).Value = "Pobreza"
$ws.Cells.Item(74, 4).Value = 25.32166776575436
$ws.Cells.Item(75, 1).Value = 2019
$ws.Cells.Item(75, 2).Value = "Costa"
$ws.Cells.Item(75, 3).Value = "Pobreza extrema"
$ws.Cells.Item(75, 4).Value = 7.179265853177669
$ws.Cells.Item(76, 1).Value = 2019
$ws.Cells.Item(76, 2).Value = "Oriente"
$ws.Cells.Item(76, 3).Value = "Pobreza"
$ws.Cells.Item(76, 4).Value = 43.81960109911133
$ws.Cells.Item(77, 1).Value = 2019
$ws.Cells.Item(77, 2).Value = "Oriente"
$ws.Cells.Item(77, 3).Value = "Pobreza extrema"
$ws.Cells.Item(77, 4).Value = 25.39897447109997
$ws.Cells.Item(78, 1).Value = 2019
$ws.Cells.Item(78, 2).Value = "Sierra"
$ws.Cells.Item(78, 3).Value = "Pobreza"
$ws.Cells.Item(78, 4).Value = 22.50522186880941
$ws.Cells.Item(79, 1).Value = 2019
$ws.Cells.Item(79, 2).Value = "Sierra"
$ws.Cells.Item(79, 3).Value = "Pobreza extrema"
$ws.Cells.Item(79, 4).Value = 8.776254627585374
$ws.Cells.Item(80, 1).Value = 2020
$ws.Cells.Item(80, 2).Value = "Costa"
$ws.Cells.Item(80, 3).Value = "Pobreza"
$ws.Cells.Item(80, 4).Value = 30.83986650715876
$ws.Cells.Item(81, 1).Value = 2020
$ws.Cells.Item(81, 2).Value = "Costa"
$ws.Cells.Item(81, 3).Value = "Pobreza extrema"
$ws.Cells.Item(81, 4).Value = 11.43275826317707
$ws.Cells.Item(82, 1).Value = 2020
$ws.Cells.Item(82, 2).Value = "Oriente"
$ws.Cells.Item(82, 3).Value = "Pobreza"
$ws.Cells.Item(82, 4).Value = 61.76635513536574
$ws.Cells.Item(83, 1).Value = 2020
$ws.Cells.Item(83, 2).Value = "Oriente"
$ws.Cells.Item(83, 3).Value = "Pobreza extrema"
$ws.Cells.Item(83, 4).Value = 44.29573755641636
$ws.Cells.Item(84, 1).Value = 2020
$ws.Cells.Item(84, 2).Value = "Sierra"
$ws.Cells.Item(84, 3).Value = "Pobreza"
$ws.Cells.Item(84, 4).Value = 25.91388848689958
$ws.Cells.Item(85, 1).Value = 2020
$ws.Cells.Item(85, 2).Value = "Sierra"
$ws.Cells.Item(85, 3).Value = "Pobreza extrema"
$ws.Cells.Item(85, 4).Value = 11.02669502542741
$ws.Cells.Item(86, 1).Value = 2021
$ws.Cells.Item(86, 2).Value = "Costa"
$ws.Cells.Item(86, 3).Value = "Pobreza"
$ws.Cells.Item(86, 4).Value = 27.30676549256789
$ws.Cells.Item(87, 1).Value = 2021
$ws.Cells.Item(87, 2).Value = "Costa"
$ws.Cells.Item(87, 3).Value = "Pobreza extrema"
$ws.Cells.Item(87, 4).Value = 7.078111881995505
$ws.Cells.Item(88, 1).Value = 2021
$ws.Cells.Item(88, 2).Value = "Oriente"
$ws.Cells.Item(88, 3).Value = "Pobreza"
$ws.Cells.Item(88, 4).Value = 59.60492083423075
$ws.Cells.Item(89, 1).Value = 2021
$ws.Cells.Item(89, 2).Value = "Oriente"
$ws.Cells.Item(89, 3).Value = "Pobreza extrema"
$ws.Cells.Item(89, 4).Value = 43.01085351732007
$ws.Cells.Item(90, 1).Value = 2021
$ws.Cells.Item(90, 2).Value = "Sierra"
$ws.Cells.Item(90, 3).Value = "Pobreza"
$ws.Cells.Item(90, 4).Value = 17.87683994527823
$ws.Cells.Item(91, 1).Value = 2021
$ws.Cells.Item(91, 2).Value = "Sierra"
$ws.Cells.Item(91, 3).Value = "Pobreza extrema"
$ws.Cells.Item(91, 4).Value = 5.123392169594305
$ws.Cells.Item(92, 1).Value = 2022
$ws.Cells.Item(92, 2).Value = "Costa"
$ws.Cells.Item(92, 3).Value = "Pobreza"
$ws.Cells.Item(92, 4).Value = 23.35297753197948
$ws.Cells.Item(93, 1).Value = 2022
$ws.Cells.Item(93, 2).Value = "Costa"
$ws.Cells.Item(93, 3).Value = "Pobreza extrema"
$ws.Cells.Item(93, 4).Value = 5.019706540264221
$ws.Cells.Item(94, 1).Value = 2022
$ws.Cells.Item(94, 2).Value = "Oriente"
$ws.Cells.Item(94, 3).Value = "Pobreza"
$ws.Cells.Item(94, 4).Value = 58.99761087461328
$ws.Cells.Item(95, 1).Value = 2022
$ws.Cells.Item(95, 2).Value = "Oriente"
$ws.Cells.Item(95, 3).Value = "Pobreza extrema"
$ws.Cells.Item(95, 4).Value = 36.6492121534727
$ws.Cells.Item(96, 1).Value = 2022
$ws.Cells.Item(96, 2).Value = "Sierra"
$ws.Cells.Item(96, 3).Value = "Pobreza"
$ws.Cells.Item(96, 4).Value = 16.86130685979331
$ws.Cells.Item(97, 1).Value = 2022
$ws.Cells.Item(97, 2).Value = "Sierra"
$ws.Cells.Item(97, 3).Value = "Pobreza extrema"
$ws.Cells.Item(97, 4).Value = 3.808267556351208
$ws.Cells.Item(98, 1).Value = 2023
$ws.Cells.Item(98, 2).Value = "Costa"
$ws.Cells.Item(98, 3).Value = "Pobreza"
$ws.Cells.Item(98, 4).Value = 23.3355753923395
$ws.Cells.Item(99, 1).Value = 2023
$ws.Cells.Item(99, 2).Value = "Costa"
$ws.Cells.Item(99, 3).Value = "Pobreza extrema"
$ws.Cells.Item(99, 4).Value = 4.967560401145677
$ws.Cells.Item(100, 1).Value = 2023
$ws.Cells.Item(100, 2).Value = "Oriente"
$ws.Cells.Item(100, 3).Value = "Pobreza"
$ws.Cells.Item(100, 4).Value = 58.56001593192315
$ws.Cells.Item(101, 1).Value = 2023
$ws.Cells.Item(101, 2).Value = "Oriente"
$ws.Cells.Item(101, 3).Value = "Pobreza extrema"
$ws.Cells.Item(101, 4).Value = 40.19148822551381
$ws.Cells.Item(102, 1).Value = 2023
$ws.Cells.Item(102, 2).Value = "Sierra"
$ws.Cells.Item(102, 3).Value = "Pobreza"
$ws.Cells.Item(102, 4).Value = 17.73067276597218
$ws.Cells.Item(103, 1).Value = 2023
$ws.Cells.Item(103, 2).Value = "Sierra"
$ws.Cells.Item(103, 3).Value = "Pobreza extrema"
$ws.Cells.Item(103, 4).Value = 5.220391528618548
$ws.Cells.Item(104, 1).Value = 2024
$ws.Cells.Item(104, 2).Value = "Costa"
$ws.Cells.Item(104, 3).Value = "Pobreza"
$ws.Cells.Item(104, 4).Value = 24.48583733280943
$ws.Cells.Item(105, 1).Value = 2024
$ws.Cells.Item(105, 2).Value = "Costa"
$ws.Cells.Item(105, 3).Value = "Pobreza extrema"
$ws.Cells.Item(105, 4).Value = 8.075678831010686
$ws.Cells.Item(106, 1).Value = 2024
$ws.Cells.Item(106, 2).Value = "Oriente"
$ws.Cells.Item(106, 3).Value = "Pobreza"
$ws.Cells.Item(106, 4).Value = 57.74246455466749
$ws.Cells.Item(107, 1).Value = 2024
$ws.Cells.Item(107, 2).Value = "Oriente"
$ws.Cells.Item(107, 3).Value = "Pobreza extrema"
$ws.Cells.Item(107, 4).Value = 40.1328197825313
$ws.Cells.Item(108, 1).Value = 2024
$ws.Cells.Item(108, 2).Value = "Sierra"
$ws.Cells.Item(108, 3).Value = "Pobreza"
$ws.Cells.Item(108, 4).Value = 19.99351893445714
$ws.Cells.Item(109, 1).Value = 2024
$ws.Cells.Item(109, 2).Value = "Sierra"
$ws.Cells.Item(109, 3).Value = "Pobreza extrema"
$ws.Cells.Item(109, 4).Value = 6.860894057451052

Write-Output "Update complete"